$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "kak"

# The phone number below looks like a numeric literal, so a plain
# Range.Value assignment would get auto-converted to a Number by Excel's
# input parser. Force text storage (as typing it in with Format Cells set
# to Text would) by flipping the cell to a text number format first, then
# clearing that formatting back off once the literal text value has been
# committed - the stored value stays text while the cell's style reverts
# to the sheet default.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "+380958630689"
$ws.Range("B1").ClearFormats()

$ws.Range("C1").Value = 0
